$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (Chair)
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 2200
$ws.Range("F2").Value = 3120
$ws.Range("G2").Value = 390
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 35
$ws.Range("K2").Value = 535

# Update row 3 (Table)
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 1660
$ws.Range("F3").Value = 2500
$ws.Range("G3").Value = 500

# Update row 4 (Chandelier)
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 405
$ws.Range("E4").Value = 405
$ws.Range("F4").Value = 625
$ws.Range("G4").Value = 625

# Remove row 5 (Lamp) entirely
$ws.Rows("5").Delete()
